$wb = $excel.ActiveWorkbook

$wsAlu = $wb.Worksheets.Item("ALU Codes")
$wsAsm = $wb.Worksheets.Item("ASSEMBLY")
$wsBus = $wb.Worksheets.Item("BUS Ctrl")

# --- ASSEMBLY sheet: new "mv" and "memlda" rows -------------------------
# Values are entered in the exact order needed so that new shared-string
# entries are created in the same sequence as in the target workbook.
$wsAsm.Range("D18").Value = "move content from A to B"
$wsAsm.Range("B18").Value = "0x31"

$wsBus.Range("A9").Value = "GP1"
$wsBus.Range("C9").Value = "010"
$wsBus.Range("D9").Value = "000010"

$wsAsm.Range("C18").Value = "mv"
$wsAsm.Range("F18").Value = "mv <op1> <op2>"
$wsAsm.Range("G18").Value = "op1 = source register (in 4 bit value), op2 = destination register (in 4 bit value)"

$wsAsm.Range("C20").Value = "memlda"
$wsAsm.Range("B20").Value = "0x35"
$wsAsm.Range("D20").Value = "load from memory to A"
$wsAsm.Range("F20").Value = "memlda <op1>"

$wsBus.Range("A10").Value = "FTREG"

# Reused shared strings (existing entries, order does not affect indices)
$wsBus.Range("B9").Value = "011"
$wsBus.Range("B10").Value = "111"
$wsBus.Range("C10").Value = "111"

# --- Column D width on ASSEMBLY sheet ------------------------------------
$wsAsm.Columns("D").ColumnWidth = 24.7

# --- Selections / active sheet -------------------------------------------
$wsAlu.Rows("13").Select()
$wsAsm.Range("F21").Select()
$wsBus.Range("J31").Select()
$wsBus.Activate()
